# Updating price policy code
# The "Electric_boiler" technology column is dropped from the Installation
# and Capacity sheets (its shared string is removed and every later column
# shifts left by one). The remaining numeric results on those two sheets,
# plus the storage-capacity figure, are refreshed with the re-run model
# output.

$wb = $excel.ActiveWorkbook

$installation = $wb.Worksheets.Item("Installation")
$capacity     = $wb.Worksheets.Item("Capacity")
$storage      = $wb.Worksheets.Item("Storage_capacity")

# Drop the "Electric_boiler" column (column B) on both sheets - this shifts
# Gas_CHP / Gas_boiler / Grid / Heat_pump / Solar_PV / Solar_thermal one
# column to the left (B:G) and removes the now-unused shared string.
$installation.Range("B1").EntireColumn.Delete()
$capacity.Range("B1").EntireColumn.Delete()

$cols = @("B", "C", "D", "E", "F", "G")

# Installation sheet: refreshed 0/1 capacity-installed flags.
$installRow2 = @(1, 0, 0, 0, 1, 0)
$installRow3 = @(1, 1, 0, 1, 0, 1)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $installation.Range("$($cols[$i])2").Value = $installRow2[$i]
    $installation.Range("$($cols[$i])3").Value = $installRow3[$i]
}

# Capacity sheet: refreshed installed-capacity results.
$capRow2 = @(49.71098265895953, 0, 0, 0, 1129.7522478085427, 0)
$capRow3 = @(86, 494.17985315754504, 0, 0, 0, 0)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $capacity.Range("$($cols[$i])2").Value = $capRow2[$i]
    $capacity.Range("$($cols[$i])3").Value = $capRow3[$i]
}

# Storage_capacity sheet: refreshed storage capacity result for Heat.
$storage.Range("B2").Value = 939.88908271631851
